# Add more (fake/sample) barcodes to each lot tab of the inventory template
# so the sheet format is easier to see at a glance.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sheet3" (2nd tab, physical file xl/worksheets/sheet2.xml)
# Extend the XE02933-series barcode list down through row 13.
# ---------------------------------------------------------------------------
$wsSheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3Vals = @("XE02933","XE02934","XE02935","XE02936","XE02937","XE02938","XE02939","XE02940","XE02941","XE02942","XE02943")
$r = 3
foreach ($v in $sheet3Vals) {
    $wsSheet3.Cells.Item($r, 2).Value = $v
    $r++
}
[void]$wsSheet3.Range("B3:B13").Select()

# ---------------------------------------------------------------------------
# Sheet "Sheet4" (3rd tab, physical file xl/worksheets/sheet3.xml)
# Replace the old single barcode with the XE30035-series list through row 16.
# ---------------------------------------------------------------------------
$wsSheet4 = $wb.Worksheets.Item("Sheet4")
$sheet4Vals = @("XE30035","XE30036","XE30037","XE30038","XE30039","XE30040","XE30041","XE30042","XE30043","XE30044","XE30045","XE30046","XE30047","XE30048")
$r = 3
foreach ($v in $sheet4Vals) {
    $wsSheet4.Cells.Item($r, 2).Value = $v
    $r++
}
[void]$wsSheet4.Range("B3:B16").Select()

# ---------------------------------------------------------------------------
# Sheet "Sheet1" (1st tab, physical file xl/worksheets/sheet1.xml)
# Replace the old single barcode with the XE67990-series list through row 21.
# ---------------------------------------------------------------------------
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1Vals = @("XE67990","XE67991","XE67992","XE67993","XE67994","XE67995","XE67996","XE67997","XE67998","XE67999","XE68000","XE68001","XE68002","XE68003","XE68004","XE68005","XE68006","XE68007","XE68008")
$r = 3
foreach ($v in $sheet1Vals) {
    $wsSheet1.Cells.Item($r, 2).Value = $v
    $r++
}

# Sheet1 is the sheet left active/selected when the workbook was saved,
# with the cursor resting on F11 (no particular data there - just where the
# cursor happened to be).
[void]$wsSheet1.Range("F11").Select()
